# Update the ZROZ historical-data sheet:
#  - row 201 volume correction
#  - row 202 High/Low/Close/Adj Close/Volume corrections
#  - two newly appended trading days (203: 2023-10-20, 204: 2023-10-23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells in row 201 and row 202 with corrected/revised data
$ws.Range("G201").Value = 506700

$ws.Range("C202").Value = 67.519997
$ws.Range("D202").Value = 64.800003
$ws.Range("E202").Value = 64.83000199999999
$ws.Range("F202").Value = 64.83000199999999
$ws.Range("G202").Value = 781500

# Add new row 203: 2023-10-20
# (force literal text so the date-looking string isn't auto-converted
#  to a date serial, then restore the default "Normal" style so no
#  leftover number-format style id is left on the cell)
$ws.Range("A203").NumberFormat = "@"
$ws.Range("A203").Value = "2023-10-20"
$ws.Range("A203").NumberFormat = "General"
$ws.Range("A203").Style = "Normal"
$ws.Range("B203").Value = 64.93000000000001
$ws.Range("C203").Value = 65.449997
$ws.Range("D203").Value = 64.540001
$ws.Range("E203").Value = 65.150002
$ws.Range("F203").Value = 65.150002
$ws.Range("G203").Value = 550100

# Add new row 204: 2023-10-23
$ws.Range("A204").NumberFormat = "@"
$ws.Range("A204").Value = "2023-10-23"
$ws.Range("A204").NumberFormat = "General"
$ws.Range("A204").Style = "Normal"
$ws.Range("B204").Value = 64.93000000000001
$ws.Range("C204").Value = 64.5
$ws.Range("D204").Value = 64.26840199999999
$ws.Range("E204").Value = 64.277496
$ws.Range("F204").Value = 64.277496
$ws.Range("G204").Value = 66718
